$d = $word.ActiveDocument

# 1. Add name to paragraph 7's empty run, then add a new empty paragraph after it.
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertBefore("Ahmed Wael Nagy Wanas - 20206008")
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()

# 2. Rename images (docPr name attribute) to match the diff's swap.
$body1 = $d.InlineShapes.Item(1)
$body1.Select() | Out-Null
$word.Selection.InlineShapes.Item(1).Name = "image3.png"

$hdr = $d.Sections.Item(1).Headers.Item(1)
$h1 = $hdr.Range.InlineShapes.Item(1)
$h1.Name = "image1.png"

$h2 = $hdr.Range.InlineShapes.Item(2)
$h2.Name = "image4.png"

Write-Output "done"
